# "Adapt tests to control version"
# Add a "version" column to the settings sheet (form_title | form_id | version)
# with value 1, and make the settings sheet the active sheet/selection.

$wb = $excel.ActiveWorkbook
$settings = $wb.Worksheets.Item("settings")

$settings.Cells.Item(1, 3).Value = "version"
$settings.Cells.Item(2, 3).Value = 1

$settings.Activate() | Out-Null
$settings.Range("C3").Select() | Out-Null
